# Updates the crypto price/volume snapshot values in columns D (Price) and
# E (Volume(1h)) to match the refreshed data pulled by the scraping job.
#
# The cells are stored as plain text (not numbers) in the workbook, so we
# temporarily force a text number format on the affected columns before
# writing the values. This stops Excel from "helpfully" re-interpreting
# strings like "308.24" or "0.12%" as numeric/percentage values. Once the
# values are written we restore the cell style back to Normal so no visible
# formatting change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceRange = $ws.Range("D2:D51")
$volumeRange = $ws.Range("E2:E51")

$priceRange.NumberFormat = "@"
$volumeRange.NumberFormat = "@"

$ws.Range("D2").Value = "308.24"
$ws.Range("E2").Value = "0.12%"
$ws.Range("D3").Value = "40.77"
$ws.Range("E3").Value = "1.84%"
$ws.Range("D4").Value = "5.117"
$ws.Range("E4").Value = "-0.26%"
$ws.Range("E5").Value = "-1.53%"
$ws.Range("D6").Value = "1.606"
$ws.Range("D7").Value = "0.9039"
$ws.Range("E7").Value = "2.62%"
$ws.Range("E8").Value = "0.29%"
$ws.Range("D9").Value = "0.1105"
$ws.Range("E9").Value = "9.27%"
$ws.Range("D10").Value = "0.1782"
$ws.Range("E10").Value = "1.93%"
$ws.Range("D11").Value = "0.09187"
$ws.Range("E11").Value = "2.48%"
$ws.Range("D12").Value = "0.04181"
$ws.Range("E12").Value = "-4.97%"
$ws.Range("D13").Value = "0.1052"
$ws.Range("E13").Value = "-0.28%"
$ws.Range("D14").Value = "0.001251"
$ws.Range("E14").Value = "-0.55%"
$ws.Range("D15").Value = "0.005835"
$ws.Range("E15").Value = "0.88%"
$ws.Range("D16").Value = "3.353"
$ws.Range("E16").Value = "-0.05%"
$ws.Range("D17").Value = "4.249"
$ws.Range("E17").Value = "0.01%"
$ws.Range("E18").Value = "-0.88%"
$ws.Range("D19").Value = "6.589"
$ws.Range("E19").Value = "-5.89%"
$ws.Range("D20").Value = "0.1364"
$ws.Range("E20").Value = "1.81%"
$ws.Range("D21").Value = "0.2760"
$ws.Range("E21").Value = "-2.90%"
$ws.Range("D22").Value = "0.04062"
$ws.Range("E22").Value = "-2.62%"
$ws.Range("D23").Value = "0.001231"
$ws.Range("E23").Value = "2.49%"
$ws.Range("D24").Value = "0.004104"
$ws.Range("E24").Value = "-0.04%"
$ws.Range("D25").Value = "0.0001301"
$ws.Range("E25").Value = "0.07%"
$ws.Range("D38").Value = "0.02422"
$ws.Range("E38").Value = "2.66%"
$ws.Range("D39").Value = "0.05184"
$ws.Range("E39").Value = "0.75%"
$ws.Range("D40").Value = "0.007766"
$ws.Range("E40").Value = "-2.29%"
$ws.Range("D41").Value = "0.1302"
$ws.Range("E41").Value = "-1.83%"
$ws.Range("D42").Value = "0.006889"
$ws.Range("E42").Value = "7.86%"
$ws.Range("D43").Value = "0.001951"
$ws.Range("E43").Value = "-1.37%"
$ws.Range("D44").Value = "0.008800"
$ws.Range("E44").Value = "-0.98%"
$ws.Range("D45").Value = "0.3329"
$ws.Range("E45").Value = "-0.11%"
$ws.Range("D46").Value = "0.00006940"
$ws.Range("E46").Value = "5.72%"
$ws.Range("E47").Value = "0.06%"
$ws.Range("D48").Value = "0.03111"
$ws.Range("E48").Value = "388.67%"
$ws.Range("D49").Value = "0.004199"
$ws.Range("E49").Value = "-40.06%"
$ws.Range("E50").Value = "0.06%"
$ws.Range("E51").Value = "0.06%"

$priceRange.Style = "Normal"
$volumeRange.Style = "Normal"
